# Auto-generated edit script applying numeric updates to the Leve profit tables
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 221.6875
$ws.Range("I103").Value = 214.14285
$ws.Range("J103").Value = 227.55556
$ws.Range("K103").Value = 642.4285500000001
$ws.Range("L103").Value = 682.66668
$ws.Range("M103").Value = -56.42855000000009
$ws.Range("N103").Value = -1854.66668

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6214356
$ws.Range("I132").Value = 7146069.5
$ws.Range("J132").Value = 2933
$ws.Range("K132").Value = 21438208.5
$ws.Range("L132").Value = 8799
$ws.Range("M132").Value = -21435678.5
$ws.Range("N132").Value = -13859

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2598.182
$ws.Range("I137").Value = 3120
$ws.Range("J137").Value = 2300
$ws.Range("K137").Value = 9360
$ws.Range("L137").Value = 6900
$ws.Range("M137").Value = -6810
$ws.Range("N137").Value = -12000

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1816.5416
$ws.Range("I138").Value = 1392.5
$ws.Range("J138").Value = 2664.625
$ws.Range("K138").Value = 4177.5
$ws.Range("L138").Value = 7993.875
$ws.Range("M138").Value = 962.5
$ws.Range("N138").Value = -18273.875

# ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 51500
$ws.Range("J139").Value = 51500
$ws.Range("L139").Value = 51500
$ws.Range("N139").Value = -61780

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23000.23
$ws.Range("I32").Value = 25633.29
$ws.Range("K32").Value = 25633.29
$ws.Range("M32").Value = -25346.29

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1908.8235
$ws.Range("I61").Value = 768.1818
$ws.Range("K61").Value = 768.1818
$ws.Range("M61").Value = -556.1818

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 733.64514
$ws.Range("I74").Value = 592.6070999999999
$ws.Range("K74").Value = 592.6070999999999
$ws.Range("M74").Value = 281.3929000000001

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 733.64514
$ws.Range("I77").Value = 592.6070999999999
$ws.Range("K77").Value = 2963.0355
$ws.Range("M77").Value = 1404.9645

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4715.516
$ws.Range("I132").Value = 5651.476
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 16954.428
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -14424.428
$ws.Range("N132").Value = -13310

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1908.8235
$ws.Range("I136").Value = 768.1818
$ws.Range("K136").Value = 2304.5454
$ws.Range("M136").Value = 245.4546

# BSM row 9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 12980
$ws.Range("J9").Value = 12980
$ws.Range("L9").Value = 12980
$ws.Range("N9").Value = -13316

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1172.3
$ws.Range("I107").Value = 1127.875
$ws.Range("J107").Value = 1350
$ws.Range("K107").Value = 1127.875
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = 792.125
$ws.Range("N107").Value = -5190

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 52576.1
$ws.Range("I134").Value = 79294
$ws.Range("J134").Value = 2957.1428
$ws.Range("K134").Value = 237882
$ws.Range("L134").Value = 8871.428400000001
$ws.Range("M134").Value = -235347
$ws.Range("N134").Value = -13941.4284

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1744.2858
$ws.Range("I16").Value = 842
$ws.Range("K16").Value = 842
$ws.Range("M16").Value = -555

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1744.2858
$ws.Range("I113").Value = 842
$ws.Range("K113").Value = 842
$ws.Range("M113").Value = 1328

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2475.0908
$ws.Range("I132").Value = 1673.7646
$ws.Range("K132").Value = 5021.293799999999
$ws.Range("M132").Value = -2491.293799999999

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1262.5294
$ws.Range("I134").Value = 1090.2142
$ws.Range("J134").Value = 2066.6667
$ws.Range("K134").Value = 3270.6426
$ws.Range("L134").Value = 6200.000100000001
$ws.Range("M134").Value = -735.6425999999997
$ws.Range("N134").Value = -11270.0001

# CUL row 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 136.14285
$ws.Range("I50").Value = 119
$ws.Range("J50").Value = 179
$ws.Range("K50").Value = 357
$ws.Range("L50").Value = 537
$ws.Range("M50").Value = 124
$ws.Range("N50").Value = -1499

# CUL row 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 136.14285
$ws.Range("I53").Value = 119
$ws.Range("J53").Value = 179
$ws.Range("K53").Value = 357
$ws.Range("L53").Value = 537
$ws.Range("M53").Value = 124
$ws.Range("N53").Value = -1499

# CUL row 54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2775
$ws.Range("I55").Value = 2350
$ws.Range("J55").Value = 2916.6667
$ws.Range("K55").Value = 7050
$ws.Range("L55").Value = 8750.000100000001
$ws.Range("M55").Value = -6873
$ws.Range("N55").Value = -9104.000100000001

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4454.5454
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 4700
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 14100
$ws.Range("M70").Value = -5685
$ws.Range("N70").Value = -14730

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 4454.5454
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 4700
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 14100
$ws.Range("M73").Value = -4908
$ws.Range("N73").Value = -16284

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4449190.5
$ws.Range("J131").Value = 7408370.5
$ws.Range("L131").Value = 22225111.5
$ws.Range("N131").Value = -22235191.5

# GSM row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19800
$ws.Range("J46").Value = 19800
$ws.Range("L46").Value = 19800
$ws.Range("N46").Value = -20112

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19800
$ws.Range("J57").Value = 19800
$ws.Range("L57").Value = 19800
$ws.Range("N57").Value = -21440

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 547.8570999999999
$ws.Range("I107").Value = 565.25
$ws.Range("K107").Value = 565.25
$ws.Range("M107").Value = 1354.75

# GSM row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8181.375
$ws.Range("I132").Value = 9841.375
$ws.Range("J132").Value = 3201.375
$ws.Range("K132").Value = 29524.125
$ws.Range("L132").Value = 9604.125
$ws.Range("M132").Value = -26994.125
$ws.Range("N132").Value = -14664.125

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6422.909
$ws.Range("I136").Value = 8450.286
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 25350.858
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -22800.858
$ws.Range("N136").Value = -13725

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1142.6666
$ws.Range("I122").Value = 1094.4445
$ws.Range("K122").Value = 3283.3335
$ws.Range("M122").Value = -833.3335000000002

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1663
$ws.Range("I132").Value = 1199.0555
$ws.Range("K132").Value = 3597.1665
$ws.Range("M132").Value = -1067.1665

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8443.177
$ws.Range("I136").Value = 10448.77
$ws.Range("K136").Value = 31346.31
$ws.Range("M136").Value = -28796.31

# WVR row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 73032
$ws.Range("J138").Value = 73032
$ws.Range("L138").Value = 73032
$ws.Range("N138").Value = -83312
